$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.548.70"
$ws.Range("E2").Value = "  +1.20%  "

$ws.Range("D3").Value = "2.532.40"
$ws.Range("E3").Value = "  +8.43%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.93"
$ws.Range("E5").Value = "  +1.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.30"
$ws.Range("E6").Value = "  +4.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  +6.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +11.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.94"
$ws.Range("E10").Value = "  +13.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0825"
$ws.Range("E11").Value = "  +4.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  +11.78%  "

$ws.Range("D13").Value = "2.917.59"
$ws.Range("E13").Value = "  +8.28%  "

$ws.Range("E14").Value = "  +2.80%  "

$ws.Range("D15").Value = "2.547.70"
$ws.Range("E15").Value = "  +8.86%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.893"
$ws.Range("E16").Value = "  +11.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.96"
$ws.Range("E17").Value = "  +9.91%  "

$ws.Range("D18").Value = "46.530.61"
$ws.Range("E18").Value = "  +1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.02"
$ws.Range("E19").Value = "  +12.37%  "

$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  +3.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.57"
$ws.Range("E21").Value = "  +10.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.62"
$ws.Range("E22").Value = "  +5.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.28"
$ws.Range("E23").Value = "  +4.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +6.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.11"
$ws.Range("E25").Value = "  +12.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.88"
$ws.Range("E27").Value = "  +3.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.43"
$ws.Range("E28").Value = "  +17.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.43"
$ws.Range("E29").Value = "  +8.64%  "

$ws.Range("E30").Value = "  +2.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.75"
$ws.Range("E31").Value = "  +3.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.98"
$ws.Range("E32").Value = "  +11.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  +4.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0841"
$ws.Range("E34").Value = "  +9.96%  "

$ws.Range("E35").Value = "  +23.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.09"
$ws.Range("E36").Value = "  +4.16%  "

$ws.Range("E37").Value = "  +7.48%  "

$ws.Range("E38").Value = "  +4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.37"
$ws.Range("E39").Value = "  +8.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.28"
$ws.Range("E40").Value = "  +11.87%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0324"
$ws.Range("E41").Value = "  +9.37%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.54"
$ws.Range("E42").Value = "  +12.12%  "

$ws.Range("D43").Value = "1.994.95"
$ws.Range("E43").Value = "  +7.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.52"
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.46"
$ws.Range("E46").Value = "  +36.06%  "

$ws.Range("E47").Value = "  +4.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.201"
$ws.Range("E48").Value = "  +9.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.42"
$ws.Range("E49").Value = "  +11.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.83"
$ws.Range("E50").Value = "  +11.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.72"
$ws.Range("E51").Value = "  +6.70%  "

# Reset style to Normal for forced-text cells to avoid extraneous style attrs
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
